$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a "Justification" style formula column (M) that concatenates a
# sentence using the tone/category value already present in column K.
# Row 3: standalone formula
$ws.Range("M3").Formula = "=CONCAT(""Use of technical language with "",K3,"" outcomes."")"

# Row 4 keeps its existing static text (no formula) - left untouched.

# Row 5: standalone formula
$ws.Range("M5").Formula = "=CONCAT(""Use of technical language with "",K5,"" outcomes."")"

# Rows 6-61: shared formula block
$ws.Range("M6:M61").Formula = "=CONCAT(""Use of technical language with "",K6,"" outcomes."")"

# Row 62 keeps its existing static text (no formula) - left untouched.

# --- Update the sheet view / selection to reflect where the author left off ---
$ws.Activate()
$ws.Range("M62").Select()
